$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: D5 (title) and E5 (link)
$ws.Range("D5").Value = "윌콕슨 순위합 검정"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/11/17/ranksum.html"

# Row 23: D23 (title) and E23 (link)
$ws.Range("D23").Value = "[터미널에서 파이썬 실행 시 파이참으로 디버깅하는 방법]How to debug code running from terminal in pycharm"
$ws.Range("E23").Value = "https://theonly1.tistory.com/2921"

# Row 29: D29 (title)
$ws.Range("D29").Value = "[만화] 인턴일기 58~65"

# Row 37: D37 (title) and E37 (link)
$ws.Range("D37").Value = "[Paper Review] Semi-Supervised Text Classification with Balanced Deep Representation Distributions"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1851&mod=document&pageid=1"
